# "added credentials to save files"
# Append three new logbook rows (3, 4, 5) to the active sheet, re-using the
# existing shared strings where the value repeats and creating new shared
# strings for genuinely new values (in the same left-to-right, top-to-bottom
# order they first appear).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Set-TextValue($cell, [string]$text) {
    # Writing a plain, non date-like string never gets reinterpreted by
    # Excel, so it can be assigned directly and lands in sharedStrings as-is.
    # Date-like strings (e.g. "2023-12-06") would otherwise be auto-converted
    # to a serial date (with an automatically-created number format style),
    # so those are entered as a literal text formula first and then replaced
    # in-place with their computed value via Copy/PasteSpecial, which yields
    # a plain shared-string cell with no style side effects.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

# Row 3: duplicate of row 2 (2023-12-06 / test123 / test123-2023-12-06.csv)
Set-TextValue $ws.Range("A3") "2023-12-06"
$ws.Range("B3").Value = "test123"
$ws.Range("C3").Value = "test123"
$ws.Range("D3").Value = "test123"
$ws.Range("E3").Value = "test123"
$ws.Range("F3").Value = "test123"
$ws.Range("G3").Value = "test123"
$ws.Range("H3").Value = "test123-2023-12-06.csv"

# Row 4: 2023-12-06 / another / another-2023-12-06.csv
Set-TextValue $ws.Range("A4") "2023-12-06"
$ws.Range("B4").Value = "another"
$ws.Range("C4").Value = "another"
$ws.Range("D4").Value = "another"
$ws.Range("E4").Value = "another"
$ws.Range("F4").Value = "another"
$ws.Range("G4").Value = "another"
$ws.Range("H4").Value = "another-2023-12-06.csv"

# Row 5: 2023-12-08 / test1 / test1-2023-12-08.csv
Set-TextValue $ws.Range("A5") "2023-12-08"
$ws.Range("B5").Value = "test1"
$ws.Range("C5").Value = "test1"
$ws.Range("D5").Value = "test1"
$ws.Range("E5").Value = "test1"
$ws.Range("F5").Value = "test1"
$ws.Range("G5").Value = "test1"
$ws.Range("H5").Value = "test1-2023-12-08.csv"
